$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format so numeric-looking strings (e.g. "514.23") are
# preserved as text, matching the original inlineStr cell contents, then
# reset the cell style back to the default "Normal" style so no extra
# formatting is left behind.
$cellRefs = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "E9", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "E32", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "B38", "C38", "D38", "E38", "B39", "C39", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.554.46'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '2.334.18'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '514.23'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '132.54'
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '0.340'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").Value = '23.61'
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '2.749.38'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '56.539.21'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("D17").Value = '2.341.64'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '10.46'
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Value = '325.69'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '4.15'
$ws.Range("E20").Value = '  -2.53%  '
$ws.Range("D21").Value = '6.72'
$ws.Range("E21").Value = '  +2.34%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '61.55'
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("D24").Value = '8.77'
$ws.Range("E24").Value = '  +11.28%  '
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("E27").Value = '  +4.41%  '
$ws.Range("D28").Value = '167.80'
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '0.0₃0722'
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '1.27'
$ws.Range("E35").Value = '  +0.85%  '
$ws.Range("D36").Value = '3.94'
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("D37").Value = '0.887'
$ws.Range("E37").Value = '  -4.48%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '1.57'
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '153.71'
$ws.Range("E39").Value = '  +12.11%  '
$ws.Range("D40").Value = '38.47'
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").Value = '280.32'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").Value = '0.0496'
$ws.Range("E46").Value = '  -1.97%  '
$ws.Range("D47").Value = '0.559'
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Value = '18.26'
$ws.Range("E48").Value = '  +5.23%  '
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("D50").Value = '0.0215'
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("D51").Value = '17.27'
$ws.Range("E51").Value = '  +2.04%  '

foreach ($ref in $cellRefs) {
    $ws.Range($ref).Style = "Normal"
}
